{"js": "// Indonesian update: Metropolis chapter re-edited.\n//\n// This script reproduces, paragraph-by-paragraph, the structural run-splits\n// introduced by the diff:\n//   1. \"Vehtari\" (ELPD paragraph) gets wrapped in <w:proofErr spellStart/spellEnd>\n//   2. \"Vehtari\" (standard-errors paragraph) gets wrapped likewise\n//   3. \"Adapt the Robinson fishing example...\" gains a red \"Done.\" remark\n//   4. \"Mathematica\" gets wrapped in <w:proofErr spellStart/spellEnd>\n//   5. \"Change transition operator...\" gains a red \"Done.\" remark\n//   6. \"...posterior space.\" is split around a relocated _GoBack bookmark\n//   7. The old _GoBack bookmark (end of doc) is removed (it moved to #6)\n//\n// Because the target shape needs precise run boundaries (and a\n// <w:proofErr/> element that has no first-class Office.js surface), each\n// touched paragraph is rebuilt in place via Range.insertOoxml(...,\n// \"Replace\") using a minimal flat-OPC wrapper. This keeps paragraph count,\n// numbering, and paragraph-level rsid/style attributes untouched while\n// giving byte-exact control over the run split.\n\nconst OOXML_HEADER =\n  `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\\n` +\n  `<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\\n` +\n  `<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\\n` +\n  `<pkg:xmlData>\\n` +\n  `<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\\n` +\n  `<w:body>\\n`;\n\nconst OOXML_FOOTER =\n  `\\n</w:body>\\n` +\n  `</w:document>\\n` +\n  `</pkg:xmlData>\\n` +\n  `</pkg:part>\\n` +\n  `</pkg:package>`;\n\nfunction wrapParagraphXml(innerParagraphXml) {\n  return OOXML_HEADER + innerParagraphXml + OOXML_FOOTER;\n}\n\nasync function replaceParagraph(paragraph, innerParagraphXml) {\n  paragraph.insertOoxml(wrapParagraphXml(innerParagraphXml), Word.InsertLocation.replace);\n  await context.sync();\n}\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// --- 1. \"Change the section on ELPD to discuss the Vehtari paper.\" ---\nawait replaceParagraph(\n  paragraphs.items[1],\n  '<w:p w:rsidR=\"00F11C84\" w:rsidRDefault=\"005F1D76\" w:rsidP=\"00F11C84\">' +\n    '<w:pPr><w:pStyle w:val=\"ListParagraph\"/><w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"1\"/></w:numPr></w:pPr>' +\n    '<w:r><w:t xml:space=\"preserve\">Change the section on ELPD to discuss the </w:t></w:r>' +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    '<w:r><w:t>Vehtari</w:t></w:r>' +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    '<w:r><w:t xml:space=\"preserve\"> paper. </w:t></w:r>' +\n    \"</w:p>\"\n);\n\n// --- 2. \"...Discuss how the standard errors work... for the Vehtari paper...\" ---\nawait replaceParagraph(\n  paragraphs.items[6],\n  '<w:p w:rsidR=\"00570388\" w:rsidRDefault=\"00EE1E19\" w:rsidP=\"00570388\">' +\n    '<w:pPr><w:pStyle w:val=\"ListParagraph\"/><w:numPr><w:ilvl w:val=\"1\"/><w:numId w:val=\"1\"/></w:numPr></w:pPr>' +\n    \"<w:r><w:t>Discuss how the standard errors work</w:t></w:r>\" +\n    '<w:r w:rsidR=\"00570388\"><w:t xml:space=\"preserve\"> (using pairwise comparisons)</w:t></w:r>' +\n    '<w:r><w:t xml:space=\"preserve\"> for the </w:t></w:r>' +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    \"<w:r><w:t>Vehtari</w:t></w:r>\" +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    '<w:r><w:t xml:space=\"preserve\"> paper, and (if it is available) how we can use methods from ANOVA to determine this optimally. </w:t></w:r>' +\n    \"</w:p>\"\n);\n\n// --- 3. \"Adapt the Robinson fishing example...\" gains a red \"Done.\" ---\nawait replaceParagraph(\n  paragraphs.items[18],\n  '<w:p w:rsidR=\"00B06904\" w:rsidRDefault=\"00B06904\" w:rsidP=\"00B06904\">' +\n    '<w:pPr><w:pStyle w:val=\"ListParagraph\"/><w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"1\"/></w:numPr></w:pPr>' +\n    \"<w:r><w:t>Adapt the Robinson fishing example to include the memory statement I included in the lectures.</w:t></w:r>\" +\n    '<w:r><w:t xml:space=\"preserve\"> </w:t></w:r>' +\n    '<w:r><w:rPr><w:color w:val=\"FF0000\"/></w:rPr><w:t>Done.</w:t></w:r>' +\n    \"</w:p>\"\n);\n\n// --- 4. \"...bees to include animations from Mathematica\" ---\nawait replaceParagraph(\n  paragraphs.items[20],\n  '<w:p w:rsidR=\"00B06904\" w:rsidRDefault=\"00B06904\" w:rsidP=\"0090153A\">' +\n    '<w:pPr><w:pStyle w:val=\"ListParagraph\"/><w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"1\"/></w:numPr></w:pPr>' +\n    '<w:r><w:t xml:space=\"preserve\">Adapt lectures\\u2019 </w:t></w:r>' +\n    '<w:r w:rsidR=\"001C3E1E\"><w:t xml:space=\"preserve\">bees to include animations from </w:t></w:r>' +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    \"<w:r><w:t>Mathematica</w:t></w:r>\" +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    \"</w:p>\"\n);\n\n// --- 5. \"Change transition operator...\" gains a red \"Done.\" ---\nawait replaceParagraph(\n  paragraphs.items[21],\n  '<w:p w:rsidR=\"001C3E1E\" w:rsidRDefault=\"001C3E1E\" w:rsidP=\"0090153A\">' +\n    '<w:pPr><w:pStyle w:val=\"ListParagraph\"/><w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"1\"/></w:numPr></w:pPr>' +\n    \"<w:r><w:t>Change transition operator to be a conditional distribution.</w:t></w:r>\" +\n    '<w:r><w:t xml:space=\"preserve\"> </w:t></w:r>' +\n    '<w:r><w:rPr><w:color w:val=\"FF0000\"/></w:rPr><w:t>Done.</w:t></w:r>' +\n    \"</w:p>\"\n);\n\n// --- 6. \"Add a figure showing NLP space vs posterior space.\" gets the\n//        _GoBack bookmark relocated inside it (mid-word split). ---\nawait replaceParagraph(\n  paragraphs.items[23],\n  '<w:p w:rsidR=\"001C36AA\" w:rsidRDefault=\"001C36AA\" w:rsidP=\"0090153A\">' +\n    '<w:pPr><w:pStyle w:val=\"ListParagraph\"/><w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"1\"/></w:numPr></w:pPr>' +\n    '<w:r><w:t xml:space=\"preserve\">Add </w:t></w:r>' +\n    '<w:r w:rsidR=\"009C4096\"><w:t>a figure showing NLP space vs posterior s</w:t></w:r>' +\n    '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/><w:bookmarkEnd w:id=\"0\"/>' +\n    \"<w:r><w:t>pace.</w:t></w:r>\" +\n    \"</w:p>\"\n);\n\n// --- 7. Remove the old _GoBack bookmark at the end of the document (it\n//        now lives in paragraph #23 above). ---\nawait replaceParagraph(\n  paragraphs.items[27],\n  '<w:p w:rsidR=\"005862DA\" w:rsidRDefault=\"005862DA\" w:rsidP=\"0090153A\">' +\n    '<w:pPr><w:pStyle w:val=\"ListParagraph\"/><w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"1\"/></w:numPr></w:pPr>' +\n    \"<w:r><w:t>Remove any integrals that can be!</w:t></w:r>\" +\n    '<w:r w:rsidR=\"00C561EA\"><w:t xml:space=\"preserve\"> Replace with sampling.</w:t></w:r>' +\n    \"</w:p>\"\n);\n", "ps1": "# Indonesian update: Metropolis chapter re-edited.\n#\n# This script reproduces, paragraph-by-paragraph, the structural run-splits\n# introduced by the diff:\n#   1. \"Vehtari\" (ELPD paragraph) gets wrapped in <w:proofErr spellStart/spellEnd>\n#   2. \"Vehtari\" (standard-errors paragraph) gets wrapped likewise\n#   3. \"Adapt the Robinson fishing example...\" gains a red \"Done.\" remark\n#   4. \"Mathematica\" gets wrapped in <w:proofErr spellStart/spellEnd>\n#   5. \"Change transition operator...\" gains a red \"Done.\" remark\n#   6. \"...posterior space.\" is split around a relocated _GoBack bookmark\n#   7. The old _GoBack bookmark (end of doc) is removed (it moved to #6)\n#\n# Because the target shape needs precise run boundaries (and a <w:proofErr/>\n# element that has no first-class COM property), each touched paragraph's\n# Range is rebuilt in place via Range.InsertXML(...) using a minimal\n# flat-OPC wrapper. InsertXML replaces the addressed range's content, so\n# calling it on Paragraphs(n).Range keeps paragraph count, numbering, and\n# paragraph-level rsid/style attributes untouched while giving byte-exact\n# control over the run split.\n\n$d = $word.ActiveDocument\n\nfunction New-FlatOpcPart([string]$paragraphXml) {\n    return @\"\n<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n<pkg:xmlData>\n<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n<w:body>\n$paragraphXml\n</w:body>\n</w:document>\n</pkg:xmlData>\n</pkg:part>\n</pkg:package>\n\"@\n}\n\n# --- 1. \"Change the section on ELPD to discuss the Vehtari paper.\" ---\n$p1 = $d.Paragraphs(2).Range\n$p1.InsertXML((New-FlatOpcPart @'\n<w:p w:rsidR=\"00F11C84\" w:rsidRDefault=\"005F1D76\" w:rsidP=\"00F11C84\"><w:pPr><w:pStyle w:val=\"ListParagraph\"/><w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"1\"/></w:numPr></w:pPr><w:r><w:t xml:space=\"preserve\">Change the section on ELPD to discuss the </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>Vehtari</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\"> paper. </w:t></w:r></w:p>\n'@))\n\n# --- 2. \"...Discuss how the standard errors work... for the Vehtari paper...\" ---\n$p2 = $d.Paragraphs(7).Range\n$p2.InsertXML((New-FlatOpcPart @'\n<w:p w:rsidR=\"00570388\" w:rsidRDefault=\"00EE1E19\" w:rsidP=\"00570388\"><w:pPr><w:pStyle w:val=\"ListParagraph\"/><w:numPr><w:ilvl w:val=\"1\"/><w:numId w:val=\"1\"/></w:numPr></w:pPr><w:r><w:t>Discuss how the standard errors work</w:t></w:r><w:r w:rsidR=\"00570388\"><w:t xml:space=\"preserve\"> (using pairwise comparisons)</w:t></w:r><w:r><w:t xml:space=\"preserve\"> for the </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>Vehtari</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\"> paper, and (if it is available) how we can use methods from ANOVA to determine this optimally. </w:t></w:r></w:p>\n'@))\n\n# --- 3. \"Adapt the Robinson fishing example...\" gains a red \"Done.\" ---\n$p3 = $d.Paragraphs(19).Range\n$p3.InsertXML((New-FlatOpcPart @'\n<w:p w:rsidR=\"00B06904\" w:rsidRDefault=\"00B06904\" w:rsidP=\"00B06904\"><w:pPr><w:pStyle w:val=\"ListParagraph\"/><w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"1\"/></w:numPr></w:pPr><w:r><w:t>Adapt the Robinson fishing example to include the memory statement I included in the lectures.</w:t></w:r><w:r><w:t xml:space=\"preserve\"> </w:t></w:r><w:r><w:rPr><w:color w:val=\"FF0000\"/></w:rPr><w:t>Done.</w:t></w:r></w:p>\n'@))\n\n# --- 4. \"...bees to include animations from Mathematica\" ---\n$p4 = $d.Paragraphs(21).Range\n$p4.InsertXML((New-FlatOpcPart @'\n<w:p w:rsidR=\"00B06904\" w:rsidRDefault=\"00B06904\" w:rsidP=\"0090153A\"><w:pPr><w:pStyle w:val=\"ListParagraph\"/><w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"1\"/></w:numPr></w:pPr><w:r><w:t xml:space=\"preserve\">Adapt lectures&#8217; </w:t></w:r><w:r w:rsidR=\"001C3E1E\"><w:t xml:space=\"preserve\">bees to include animations from </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>Mathematica</w:t></w:r><w:proofErr w:type=\"spellEnd\"/></w:p>\n'@))\n\n# --- 5. \"Change transition operator...\" gains a red \"Done.\" ---\n$p5 = $d.Paragraphs(22).Range\n$p5.InsertXML((New-FlatOpcPart @'\n<w:p w:rsidR=\"001C3E1E\" w:rsidRDefault=\"001C3E1E\" w:rsidP=\"0090153A\"><w:pPr><w:pStyle w:val=\"ListParagraph\"/><w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"1\"/></w:numPr></w:pPr><w:r><w:t>Change transition operator to be a conditional distribution.</w:t></w:r><w:r><w:t xml:space=\"preserve\"> </w:t></w:r><w:r><w:rPr><w:color w:val=\"FF0000\"/></w:rPr><w:t>Done.</w:t></w:r></w:p>\n'@))\n\n# --- 6. \"Add a figure showing NLP space vs posterior space.\" gets the\n#        _GoBack bookmark relocated inside it (mid-word split). ---\n$p6 = $d.Paragraphs(24).Range\n$p6.InsertXML((New-FlatOpcPart @'\n<w:p w:rsidR=\"001C36AA\" w:rsidRDefault=\"001C36AA\" w:rsidP=\"0090153A\"><w:pPr><w:pStyle w:val=\"ListParagraph\"/><w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"1\"/></w:numPr></w:pPr><w:r><w:t xml:space=\"preserve\">Add </w:t></w:r><w:r w:rsidR=\"009C4096\"><w:t>a figure showing NLP space vs posterior s</w:t></w:r><w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/><w:bookmarkEnd w:id=\"0\"/><w:r><w:t>pace.</w:t></w:r></w:p>\n'@))\n\n# --- 7. Remove the old _GoBack bookmark at the end of the document (it now\n#        lives in paragraph #24 above). ---\n$p7 = $d.Paragraphs(28).Range\n$p7.InsertXML((New-FlatOpcPart @'\n<w:p w:rsidR=\"005862DA\" w:rsidRDefault=\"005862DA\" w:rsidP=\"0090153A\"><w:pPr><w:pStyle w:val=\"ListParagraph\"/><w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"1\"/></w:numPr></w:pPr><w:r><w:t>Remove any integrals that can be!</w:t></w:r><w:r w:rsidR=\"00C561EA\"><w:t xml:space=\"preserve\"> Replace with sampling.</w:t></w:r></w:p>\n'@))\n"}
